$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used as a scratch area to force Excel to treat numeric-looking
# text (e.g. "1.00", "6.44") as literal text rather than coercing it to a
# number. A formula that evaluates to a string, copied and pasted as values
# only, preserves the text exactly (incl. trailing zeros) without touching
# the destination cells style/number format.
$helper = $ws.Range("ZZ1")

# Row 2
$ws.Range("D2").Value = "60.226.90"
$ws.Range("E2").Value = "  +5.47%  "

# Row 3
$ws.Range("D3").Value = "2.597.57"
$ws.Range("E3").Value = "  +7.23%  "

# Row 4
$helper.Formula = "=""1.00"""
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$helper.Formula = "=""507.31"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E5").Value = "  +3.78%  "

# Row 6
$helper.Formula = "=""156.22"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E6").Value = "  +1.58%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$helper.Formula = "=""0.589"""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E8").Value = "  -4.39%  "

# Row 9
$ws.Range("D9").Value = "2.628.12"
$ws.Range("E9").Value = "  +7.64%  "

# Row 10
$helper.Formula = "=""6.44"""
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E10").Value = "  +4.35%  "

# Row 11
$ws.Range("E11").Value = "  +4.19%  "

# Row 12
$ws.Range("E12").Value = "  +2.69%  "

# Row 13
$ws.Range("E13").Value = "  +0.90%  "

# Row 14
$ws.Range("D14").Value = "3.045.98"
$ws.Range("E14").Value = "  +7.26%  "

# Row 15
$ws.Range("D15").Value = "60.329.70"
$ws.Range("E15").Value = "  +5.58%  "

# Row 16
$helper.Formula = "=""21.65"""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E16").Value = "  +5.08%  "

# Row 17
$helper.Formula = "=""0.0000140"""
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E17").Value = "  +4.80%  "

# Row 18
$ws.Range("D18").Value = "2.621.32"
$ws.Range("E18").Value = "  +7.79%  "

# Row 19
$helper.Formula = "=""4.76"""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E19").Value = "  +2.32%  "

# Row 20
$helper.Formula = "=""343.72"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E20").Value = "  +5.74%  "

# Row 21
$helper.Formula = "=""10.44"""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E21").Value = "  +4.34%  "

# Row 22
$helper.Formula = "=""6.16"""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E22").Value = "  +3.52%  "

# Row 23
$helper.Formula = "=""0.998"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$helper.Formula = "=""60.34"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E24").Value = "  +4.35%  "

# Row 25
$helper.Formula = "=""0.423"""
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E25").Value = "  +5.02%  "

# Row 26
$ws.Range("D26").Value = "2.717.71"
$ws.Range("E26").Value = "  +7.46%  "

# Row 27
$ws.Range("E27").Value = "  +2.79%  "

# Row 28
$helper.Formula = "=""0.992"""
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E28").Value = "  -0.72%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0856"
$ws.Range("E29").Value = "  +9.06%  "

# Row 30
$helper.Formula = "=""7.55"""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E30").Value = "  +3.59%  "

# Row 31
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$helper.Formula = "=""19.47"""
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E32").Value = "  +4.40%  "

# Row 33
$helper.Formula = "=""156.22"""
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E33").Value = "  +3.09%  "

# Row 34
$ws.Range("E34").Value = "  +3.32%  "

# Row 35
$helper.Formula = "=""5.74"""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E35").Value = "  +8.25%  "

# Row 36
$helper.Formula = "=""4.02"""
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E36").Value = "  +6.44%  "

# Row 37
$ws.Range("E37").Value = "  +4.39%  "

# Row 38
$helper.Formula = "=""310.50"""
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E38").Value = "  +8.65%  "

# Row 39
$ws.Range("E39").Value = "  +7.69%  "

# Row 40
$helper.Formula = "=""0.847"""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E40").Value = "  +3.32%  "

# Row 41
$helper.Formula = "=""3.77"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E41").Value = "  +7.14%  "

# Row 42
$helper.Formula = "=""0.837"""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E42").Value = "  +27.87%  "

# Row 43
$helper.Formula = "=""35.48"""
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E43").Value = "  +4.40%  "

# Row 44
$helper.Formula = "=""0.628"""
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E44").Value = "  +5.02%  "

# Row 45
$helper.Formula = "=""0.0571"""
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E45").Value = "  +7.55%  "

# Row 46
$ws.Range("E46").Value = "  -1.35%  "

# Row 47
$helper.Formula = "=""0.992"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E47").Value = "  -0.24%  "

# Row 48
$helper.Formula = "=""19.84"""
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Clear()
$ws.Range("E48").Value = "  +12.55%  "

# Row 49
$ws.Range("E49").Value = "  +6.52%  "

# Row 50
$ws.Range("E50").Value = "  +3.48%  "

# Row 51
$ws.Range("D51").Value = "2.045.48"
$ws.Range("E51").Value = "  +7.06%  "

$ws.Application.CutCopyMode = $false